# Rodada 6 - Fase 2
# Atualização da Rodada 6 e ajuste da Fase 2 da Libertadores.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the scores (Placar Casa / Placar Visitante) for Rodada 6 (rows 52-61) ---
$scores = @{
    52 = @(3, 1)
    53 = @(2, 2)
    54 = @(1, 1)
    55 = @(0, 0)
    56 = @(2, 0)
    57 = @(4, 0)
    58 = @(0, 1)
    59 = @(1, 0)
    60 = @(1, 1)
    61 = @(1, 2)
}

foreach ($row in $scores.Keys) {
    $vals = $scores[$row]
    $ws.Range("F$row").Value = $vals[0]
    $ws.Range("G$row").Value = $vals[1]
}

# --- Append the new Rodada 7 matches (rows 62-71) ---
$newMatches = @(
    @(7, "2025-05-02", "21:30", "Morumbis", "SAO", "FOR"),
    @(7, "2025-05-03", "18:30", "Maracanã", "FLU", "SPT"),
    @(7, "2025-05-03", "18:30", "Neo Química Arena", "COR", "INT"),
    @(7, "2025-05-03", "18:30", "Presidente Vargas (CE)", "CEA", "VIT"),
    @(7, "2025-05-03", "21:00", "Casa de Apostas Arena Fonte Nova", "BAH", "BOT"),
    @(7, "2025-05-04", "16:00", "Mané Garrincha", "VAS", "PAL"),
    @(7, "2025-05-04", "16:00", "Arena do Grêmio", "GRE", "SAN"),
    @(7, "2025-05-04", "18:30", "Mineirão", "CRU", "FLA"),
    @(7, "2025-05-05", "19:00", "Cícero de Souza Marques", "RBB", "MIR"),
    @(7, "2025-05-05", "20:00", "Alfredo Jaconi", "JUV", "CAM")
)

$startRow = 62
$endRow = $startRow + $newMatches.Count - 1

# Force the "Data" column to be treated as plain text so strings like
# "2025-05-02" are not auto-converted into date serial numbers.
$ws.Range("B$startRow`:B$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newMatches.Count; $i++) {
    $r = $startRow + $i
    $m = $newMatches[$i]
    $ws.Range("A$r").Value = $m[0]
    $ws.Range("B$r").Value = $m[1]
    $ws.Range("C$r").Value = $m[2]
    $ws.Range("D$r").Value = $m[3]
    $ws.Range("E$r").Value = $m[4]
    $ws.Range("H$r").Value = $m[5]

    # Scores aren't known yet for these upcoming matches, but still
    # materialize the (blank) Placar Casa/Visitante cells, matching the
    # other rows in the sheet that carry an explicit empty cell.
    $ws.Range("F$r").Style = "Normal"
    $ws.Range("G$r").Style = "Normal"
}

# Remove the temporary text formatting again so the new cells keep the
# same (default) style as the rest of the sheet.
$ws.Range("B$startRow`:B$endRow").ClearFormats()
